# "Generate Report for handoff"
#
# The localization handoff for 27b614dc-0e40-4ea4-b7cd-308a7b2554ce.md failed,
# so a brand-new report is generated for it (new guid-named handoff file) and,
# because the handoff failed, the previously recorded handoff artifacts
# (handoff file, handoff datetime, handoff reason) on the per-language sheets
# are reset / cleared.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$oldMd    = "27b614dc-0e40-4ea4-b7cd-308a7b2554ce.md"
$newMd    = "1f5eeef4-23f6-43cc-b531-a6f094206bcb.md"
$zeroDate = "0001-01-01 00:00:00"

# --- Rename the source file everywhere (cell text + hyperlink display text) ---
$ws1.Range("A2").Value = $newMd
$ws2.Range("A2").Value = $newMd
$ws3.Range("A2").Value = $newMd

foreach ($ws in @($ws1, $ws2, $ws3)) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.TextToDisplay -eq $oldMd) {
            $hl.TextToDisplay = $newMd
        }
    }
}

# --- Status changes from "Not yet handed off" to "Handoff failed" everywhere ---
$ws1.Range("B2").Value = "Handoff failed"
$ws1.Range("C2").Value = "Handoff failed"
$ws2.Range("B2").Value = "Handoff failed"
$ws3.Range("B2").Value = "Handoff failed"

# --- Because the handoff failed, drop the stale "Latest Handoff File" link/value ---
foreach ($hl in $ws2.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$2') {
        $hl.Delete()
    }
}
foreach ($hl in $ws3.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$2') {
        $hl.Delete()
    }
}
$ws2.Range("C2").Clear()
$ws3.Range("C2").Clear()

# --- Reset the handoff/handback datetimes and the handoff reason on both rows ---
$ws2.Range("D2").Value = $zeroDate
$ws2.Range("G2").Value = $zeroDate
$ws2.Range("H2").Value = "Ignored"
$ws2.Range("D3").Value = $zeroDate
$ws2.Range("G3").Value = $zeroDate
$ws2.Range("H3").Value = "Ignored"

$ws3.Range("D2").Value = $zeroDate
$ws3.Range("G2").Value = $zeroDate
$ws3.Range("H2").Value = "Ignored"
$ws3.Range("D3").Value = $zeroDate
$ws3.Range("G3").Value = $zeroDate
$ws3.Range("H3").Value = "Ignored"
